# Adds "Color" (column D) values for the supplementary macroscopic thin-section
# table, matching the commit "Corrected draft - Added photos used for plates".
#
# The shared-string table is order-sensitive (Excel appends a new unique string
# the first time it is encountered), so the cells below are populated in the
# same first-seen order as the authoritative edit rather than strict row order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D ("Color") values, in first-use order for the shared string table ---
$ws.Range("D6").Value  = "10YR 4/1, 10 YR 5/2, 10YR 6/2"
$ws.Range("D8").Value  = "5YR 6/2"
$ws.Range("D2").Value  = "10YR 6/1, 10YR 7/1, 10YR 8/1"
$ws.Range("D3").Value  = "2.5YR 6/1, 10R 4/3, 10R 8/2"
$ws.Range("D19").Value = "10YR 5/2, 10YR 5/3"
$ws.Range("D18").Value = "10YR 7/2, 10YR 6/2"
$ws.Range("D16").Value = "10YR 8/2, 10YR 7/1"
$ws.Range("D9").Value  = "2.5Y 6/2, 2.5Y 7/4, 2.5Y 8/1, 5RP 6/2, 5RP 4/2"
$ws.Range("D10").Value = "10YR 7/1, 10YR 8/1, 10YR 8/2"
$ws.Range("D21").Value = "10YR 5/1, 10YR 6/1, 10YR 7/1"
$ws.Range("D17").Value = "10YR 8/2, 10YR 6/3"

$ws.Range("D4").Value  = "10YR 6/1, 10YR 7/1, 10YR 8/1"
$ws.Range("D5").Value  = "10YR 6/1, 10YR 7/1, 10YR 8/1"
$ws.Range("D20").Value = "10YR 5/2, 10YR 5/3"
$ws.Range("D53").Value = "N/A"

# --- A handful of rows (16, 17, 18, 53) use the workbook's secondary "Segoe UI"
#     style (same style already used by D11/D13/D14/D15/... ) instead of the
#     default font. Copy/PasteSpecial-Formats the existing style from D11 onto
#     them (rather than rebuilding the font piecemeal) so no redundant
#     font/cellXf entries get left behind in styles.xml. ---
$ws.Range("D11").Copy()
foreach ($addr in @("D16","D17","D18","D53")) {
    $ws.Range($addr).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

$ws.Rows.Item(16).RowHeight = 15
$ws.Rows.Item(17).RowHeight = 15
$ws.Rows.Item(18).RowHeight = 15
$ws.Rows.Item(53).RowHeight = 15

# --- View state: the saved selection moved from B48 to D52, scrolled so row 42
#     is the first visible row. ---
$ws.Activate()
$ws.Range("D52").Select()
$excel.ActiveWindow.ScrollRow = 42
$excel.ActiveWindow.ScrollColumn = 1
